# Update "想去人数" (F column) counts across the sheets of the
# 广州-漫展信息 workbook, reflecting refreshed scrape data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 278
$ws1.Cells.Item(4, 6).Value = 1146
$ws1.Cells.Item(5, 6).Value = 9
$ws1.Cells.Item(6, 6).Value = 2770
$ws1.Cells.Item(8, 6).Value = 703
$ws1.Cells.Item(9, 6).Value = 91
$ws1.Cells.Item(10, 6).Value = 283
$ws1.Cells.Item(11, 6).Value = 198
$ws1.Cells.Item(12, 6).Value = 707
$ws1.Cells.Item(14, 6).Value = 133
$ws1.Cells.Item(15, 6).Value = 1709
$ws1.Cells.Item(16, 6).Value = 310
$ws1.Cells.Item(18, 6).Value = 204

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(4, 6).Value = 28
$ws2.Cells.Item(5, 6).Value = 14
$ws2.Cells.Item(6, 6).Value = 17
$ws2.Cells.Item(9, 6).Value = 122
$ws2.Cells.Item(10, 6).Value = 34
$ws2.Cells.Item(12, 6).Value = 53
$ws2.Cells.Item(18, 6).Value = 35
$ws2.Cells.Item(23, 6).Value = 24

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 796
$ws3.Cells.Item(4, 6).Value = 2028
$ws3.Cells.Item(5, 6).Value = 264

# --- Sheet "全部类型" (All types, a merged view of the sheets above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 796
$ws4.Cells.Item(4, 6).Value = 2028
$ws4.Cells.Item(5, 6).Value = 264
$ws4.Cells.Item(9, 6).Value = 28
$ws4.Cells.Item(10, 6).Value = 14
$ws4.Cells.Item(11, 6).Value = 278
$ws4.Cells.Item(12, 6).Value = 1146
$ws4.Cells.Item(13, 6).Value = 9
$ws4.Cells.Item(14, 6).Value = 17
$ws4.Cells.Item(17, 6).Value = 2770
$ws4.Cells.Item(18, 6).Value = 122
$ws4.Cells.Item(20, 6).Value = 34
$ws4.Cells.Item(22, 6).Value = 53
$ws4.Cells.Item(23, 6).Value = 703
$ws4.Cells.Item(24, 6).Value = 91
$ws4.Cells.Item(25, 6).Value = 283
$ws4.Cells.Item(27, 6).Value = 198
$ws4.Cells.Item(28, 6).Value = 707
$ws4.Cells.Item(30, 6).Value = 133
$ws4.Cells.Item(32, 6).Value = 1709
$ws4.Cells.Item(33, 6).Value = 310
$ws4.Cells.Item(37, 6).Value = 204
$ws4.Cells.Item(39, 6).Value = 35
$ws4.Cells.Item(44, 6).Value = 24
